$d = $word.ActiveDocument

# 1) "summarized:" -> "summarised:"
$d.Content.Find.Execute("summarized:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "summarised:", 2) | Out-Null

# 2) "customized version of EC service (" -> "customised version of EC service ("
$d.Content.Find.Execute("customized version of EC service (", $true, $false, $false, $false, $false, `
    $true, 1, $false, "customised version of EC service (", 2) | Out-Null

# 3) "organized" -> "organised" (avoid touching the nearby apostrophe in "OP's")
$d.Content.Find.Execute("organized", $true, $false, $false, $false, $false, `
    $true, 1, $false, "organised", 2) | Out-Null

Write-Host "Done"
